$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data: the two distinct names in column A ("Vinh Duong" / "Foden
# Duong") are replaced by a single shared value "Random" for both data rows.
$ws.Range("A2").Value = "Random"
$ws.Range("A3").Value = "Random"

# Columns A and B get an explicit ("best fit") width, sized to the longest
# entry now present in each column.
$ws.Columns.Item(1).ColumnWidth = 11.83
$ws.Columns.Item(2).ColumnWidth = 11

# Move/restore the active selection to F13 (was U13).
$ws.Range("F13").Select()
